$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 98.5
$ws.Range("I4").Value = 98.5
$ws.Range("K4").Value = 98.5
$ws.Range("M4").Value = 15.5
$ws.Range("H28").Value = 14447.0625
$ws.Range("I28").Value = 5796.2144
$ws.Range("J28").Value = 75003
$ws.Range("K28").Value = 5796.2144
$ws.Range("L28").Value = 75003
$ws.Range("M28").Value = -5311.2144
$ws.Range("N28").Value = -75973
$ws.Range("H53").Value = 242.625
$ws.Range("I53").Value = 250.2
$ws.Range("J53").Value = 230
$ws.Range("K53").Value = 250.2
$ws.Range("L53").Value = 230
$ws.Range("M53").Value = 386.8
$ws.Range("N53").Value = -1504
$ws.Range("H137").Value = 2198.7378
$ws.Range("I137").Value = 2011.6222
$ws.Range("K137").Value = 6034.8666
$ws.Range("M137").Value = -3484.8666
$ws.Range("H138").Value = 2073.0564
$ws.Range("I138").Value = 1415.4
$ws.Range("J138").Value = 2921.6453
$ws.Range("K138").Value = 4246.200000000001
$ws.Range("L138").Value = 8764.9359
$ws.Range("M138").Value = 893.7999999999993
$ws.Range("N138").Value = -19044.9359

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -192
$ws.Range("N29").Value = $null
$ws.Range("H32").Value = 694878.1
$ws.Range("I32").Value = 759644.3
$ws.Range("K32").Value = 759644.3
$ws.Range("M32").Value = -759357.3
$ws.Range("H74").Value = 1548.3208
$ws.Range("I74").Value = 1248.3158
$ws.Range("K74").Value = 1248.3158
$ws.Range("M74").Value = -374.3158000000001
$ws.Range("H77").Value = 1548.3208
$ws.Range("I77").Value = 1248.3158
$ws.Range("K77").Value = 6241.579000000001
$ws.Range("M77").Value = -1873.579000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 22721.285
$ws.Range("I102").Value = 9841.5
$ws.Range("J102").Value = 100000
$ws.Range("K102").Value = 9841.5
$ws.Range("L102").Value = 100000
$ws.Range("M102").Value = -6596.5
$ws.Range("N102").Value = -106490

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5425.613
$ws.Range("I31").Value = 1100.5518
$ws.Range("J31").Value = 9226.424000000001
$ws.Range("K31").Value = 1100.5518
$ws.Range("L31").Value = 9226.424000000001
$ws.Range("M31").Value = -805.5518
$ws.Range("N31").Value = -9816.424000000001
$ws.Range("H34").Value = 5425.613
$ws.Range("I34").Value = 1100.5518
$ws.Range("J34").Value = 9226.424000000001
$ws.Range("K34").Value = 1100.5518
$ws.Range("L34").Value = 9226.424000000001
$ws.Range("M34").Value = -898.5518
$ws.Range("N34").Value = -9630.424000000001
$ws.Range("H58").Value = 1770.381
$ws.Range("I58").Value = 1292.1818
$ws.Range("J58").Value = 2296.4
$ws.Range("K58").Value = 1292.1818
$ws.Range("L58").Value = 2296.4
$ws.Range("M58").Value = -1089.1818
$ws.Range("N58").Value = -2702.4
$ws.Range("H132").Value = 3547529.2
$ws.Range("I132").Value = 1205.4722
$ws.Range("J132").Value = 15153680
$ws.Range("K132").Value = 3616.4166
$ws.Range("L132").Value = 45461040
$ws.Range("M132").Value = -1086.4166
$ws.Range("N132").Value = -45466100
$ws.Range("H134").Value = 5557.8887
$ws.Range("I134").Value = 5492.5454
$ws.Range("K134").Value = 16477.6362
$ws.Range("M134").Value = -13942.6362
$ws.Range("H135").Value = 51600
$ws.Range("I135").Value = 50000
$ws.Range("J135").Value = 53200
$ws.Range("K135").Value = 50000
$ws.Range("L135").Value = 53200
$ws.Range("M135").Value = -44930
$ws.Range("N135").Value = -63340
$ws.Range("H136").Value = 1770.381
$ws.Range("I136").Value = 1292.1818
$ws.Range("J136").Value = 2296.4
$ws.Range("K136").Value = 3876.5454
$ws.Range("L136").Value = 6889.200000000001
$ws.Range("M136").Value = -1326.5454
$ws.Range("N136").Value = -11989.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 30.777779
$ws.Range("I2").Value = 10.125
$ws.Range("J2").Value = 47.3
$ws.Range("K2").Value = 60.75
$ws.Range("L2").Value = 283.8
$ws.Range("M2").Value = 52.25
$ws.Range("N2").Value = -509.8
$ws.Range("H94").Value = 2341.182
$ws.Range("I94").Value = 813.25
$ws.Range("J94").Value = 3214.2856
$ws.Range("K94").Value = 2439.75
$ws.Range("L94").Value = 9642.856800000001
$ws.Range("M94").Value = -1763.75
$ws.Range("N94").Value = -10994.8568
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null
$ws.Range("H103").Value = 505
$ws.Range("I103").Value = 505
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 1515
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -636
$ws.Range("N103").Value = $null
$ws.Range("H106").Value = 8973.923000000001
$ws.Range("J106").Value = 8973.923000000001
$ws.Range("L106").Value = 26921.769
$ws.Range("N106").Value = -28813.769
$ws.Range("H112").Value = 6080
$ws.Range("I112").Value = 7450
$ws.Range("J112").Value = 5166.6665
$ws.Range("K112").Value = 22350
$ws.Range("L112").Value = 15499.9995
$ws.Range("M112").Value = -21242
$ws.Range("N112").Value = -17715.9995
$ws.Range("H121").Value = 1250.6666
$ws.Range("I121").Value = 435.57144
$ws.Range("J121").Value = 1535.95
$ws.Range("K121").Value = 1306.71432
$ws.Range("L121").Value = 4607.85
$ws.Range("M121").Value = 3.285679999999957
$ws.Range("N121").Value = -7227.85
$ws.Range("H122").Value = 5908.95
$ws.Range("I122").Value = 558.5
$ws.Range("J122").Value = 11259.4
$ws.Range("K122").Value = 5026.5
$ws.Range("L122").Value = 101334.6
$ws.Range("M122").Value = -2576.5
$ws.Range("N122").Value = -106234.6
$ws.Range("H125").Value = 2650.182
$ws.Range("J125").Value = 2815.2
$ws.Range("L125").Value = 8445.599999999999
$ws.Range("N125").Value = -18285.6
$ws.Range("H131").Value = 1071.5217
$ws.Range("J131").Value = 1097.5
$ws.Range("L131").Value = 3292.5
$ws.Range("N131").Value = -13372.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 54339.332
$ws.Range("I22").Value = 3000
$ws.Range("K22").Value = 3000
$ws.Range("M22").Value = -2471
$ws.Range("H126").Value = 2317.7144
$ws.Range("I126").Value = 2600
$ws.Range("J126").Value = 2270.6667
$ws.Range("K126").Value = 7800
$ws.Range("L126").Value = 6812.000100000001
$ws.Range("M126").Value = -5330
$ws.Range("N126").Value = -11752.0001
$ws.Range("H132").Value = 2676.698
$ws.Range("I132").Value = 2448.2432
$ws.Range("J132").Value = 3205
$ws.Range("K132").Value = 7344.7296
$ws.Range("L132").Value = 9615
$ws.Range("M132").Value = -4814.7296
$ws.Range("N132").Value = -14675

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 3500
$ws.Range("I107").Value = 3500
$ws.Range("K107").Value = 3500
$ws.Range("M107").Value = -1580
$ws.Range("H136").Value = 8335368
$ws.Range("I136").Value = 1891.6364
$ws.Range("J136").Value = 18520728
$ws.Range("K136").Value = 5674.9092
$ws.Range("L136").Value = 55562184
$ws.Range("M136").Value = -3124.9092
$ws.Range("N136").Value = -55567284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 356997
$ws.Range("I29").Value = 980
$ws.Range("J29").Value = 535005.5
$ws.Range("K29").Value = 980
$ws.Range("L29").Value = 535005.5
$ws.Range("M29").Value = -690
$ws.Range("N29").Value = -535585.5
$ws.Range("H95").Value = 84114.664
$ws.Range("J95").Value = 84114.664
$ws.Range("L95").Value = 84114.664
$ws.Range("N95").Value = -89606.664
$ws.Range("H126").Value = 1749.6666
$ws.Range("I126").Value = 1124.5
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 3373.5
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -903.5
$ws.Range("N126").Value = -13940
$ws.Range("H136").Value = 2409.0637
$ws.Range("I136").Value = 1883.1025
$ws.Range("J136").Value = 4973.125
$ws.Range("K136").Value = 5649.3075
$ws.Range("L136").Value = 14919.375
$ws.Range("M136").Value = -3099.3075
$ws.Range("N136").Value = -20019.375

Write-Output "Applied all changes"